# "finishing the scrum report" -- fill in the two remaining "Problems/Concerns"
# comments for the "reduce" method rows, grow those rows to fit the wrapped
# text, and leave the selection where the author's cursor ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Write G14's comment before G12's so the workbook's shared-string table ends
# up with the same two new entries in the same order (index 32, then 33).
$ws.Range("G14").Value = "the reduce function originally worled differsent than it dos not because I realized that it would be easier to just take the letters and numbers I wanted instead of take out a few that were specified."
$ws.Range("G12").Value = "reduce makes all of the letters lowercase to make the check funtion easier to use."

# Grow rows 12 and 14 so the newly-added, wrapped comment text is fully visible.
$ws.Rows.Item(12).RowHeight = 165
$ws.Rows.Item(14).RowHeight = 75

# Scroll the window down and leave the selection on I13, matching where the
# author was working when the report was saved.
$excel.ActiveWindow.ScrollRow = 7
$ws.Range("I13").Select()
